# Add "Up"/"Down" style sentiment category and a new trade row to the
# GILD sentiment sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns for row 3 (existing last data row) ---------------------
# Column X: numeric delta; Column Y: text category ("Up")
$ws.Cells.Item(3, 24).Value = 0.06999999999999318   # X3
$ws.Cells.Item(3, 25).Value = "Up"                  # Y3

# --- New row 4 with a fresh trade record --------------------------------
$ws.Cells.Item(4, 1).Value  = 42633.884317129632   # A4 Date
$ws.Cells.Item(4, 2).Value  = 3                    # B4 ScoreFinal
$ws.Cells.Item(4, 3).Value  = "Neutral"             # C4 Verdict
$ws.Cells.Item(4, 4).Value  = 24                   # D4 totalSentiment
$ws.Cells.Item(4, 5).Value  = 20499                # E4 wordCount
$ws.Cells.Item(4, 6).Value  = 980                  # F4 sentenceCount
$ws.Cells.Item(4, 7).Value  = 63                   # G4 posWordPercentage
$ws.Cells.Item(4, 8).Value  = 34                   # H4 negWordPercentage
$ws.Cells.Item(4, 9).Value  = 88                   # I4 posPhrasePercentage
$ws.Cells.Item(4, 10).Value = 11                   # J4 negPhrasePercentage
$ws.Cells.Item(4, 11).Value = 23028                # K4 ElapsedMs
$ws.Cells.Item(4, 12).Value = 219                  # L4 posWordCount
$ws.Cells.Item(4, 13).Value = 118                  # M4 negWordCount
$ws.Cells.Item(4, 14).Value = 15                   # N4 positivePhraseCount
$ws.Cells.Item(4, 15).Value = 2                    # O4 negativePhraseCount
$ws.Cells.Item(4, 16).Value = "Named"               # P4 Method
$ws.Cells.Item(4, 17).Value = 0                    # Q4 RSI
$ws.Cells.Item(4, 18).Value = -31.95                # R4 PEG
$ws.Cells.Item(4, 19).Value = -0.082                # S4 200Moving%
$ws.Cells.Item(4, 20).Value = -0.28                 # T4 50Moving%
$ws.Cells.Item(4, 21).Value = 6.77                  # U4 PriceBook
$ws.Cells.Item(4, 22).Value = 1.88                  # V4 Dividend
$ws.Cells.Item(4, 23).Value = 0                     # W4 Bollinger

# The "200Moving%" column is formatted as a percentage in every other row
# (style index 2 / numFmtId 10); give the new row the same format. Column A
# already inherits the date/time display from the column-level style, so it
# needs no extra formatting.
$ws.Cells.Item(4, 19).NumberFormat = "0.00%"
